# The source data table had a duplicate/near-duplicate "domestic gross ($million)"
# row (row 11) that is being removed as part of a clustering / "useless columns"
# cleanup pass. Deleting the entire row shifts the remaining rows (12-19) up by
# one, which matches the new dimension A1:B18.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").EntireRow.Delete()
